# Fixes #67 - Updated related resource implementation for Batch Import.
# Adds a new "related resource" column (T) to the import test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, styled the same as the other bold header cells (e.g. C1, D1 ...)
$ws.Range("T1").Value = "related resource"
$ws.Range("T1").Style = $ws.Range("C1").Style

# New data cell for the first data row
$ws.Range("T2").Value = "RELATED RESOURCE @{relatedType=relation; url=http://test.com/related_resource/relation}"

# Move the viewport / selection to match the edited area
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("T2").Select()
